$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 53751.06
$ws.Range("J17").Value = 53751.06
$ws.Range("L17").Value = 161253.18
$ws.Range("N17").Value = -161589.18
$ws.Range("H62").Value = 18777.812
$ws.Range("I62").Value = 14744.5
$ws.Range("J62").Value = 25500
$ws.Range("K62").Value = 14744.5
$ws.Range("L62").Value = 25500
$ws.Range("M62").Value = -14120.5
$ws.Range("N62").Value = -26748
$ws.Range("H65").Value = 18777.812
$ws.Range("I65").Value = 14744.5
$ws.Range("J65").Value = 25500
$ws.Range("K65").Value = 73722.5
$ws.Range("L65").Value = 127500
$ws.Range("M65").Value = -70602.5
$ws.Range("N65").Value = -133740
$ws.Range("H70").Value = 15450.223
$ws.Range("I70").Value = 20769.8
$ws.Range("J70").Value = 8800.75
$ws.Range("K70").Value = 62309.39999999999
$ws.Range("L70").Value = 26402.25
$ws.Range("M70").Value = -62039.39999999999
$ws.Range("N70").Value = -26942.25
$ws.Range("H73").Value = 15450.223
$ws.Range("I73").Value = 20769.8
$ws.Range("J73").Value = 8800.75
$ws.Range("K73").Value = 62309.39999999999
$ws.Range("L73").Value = 26402.25
$ws.Range("M73").Value = -61373.39999999999
$ws.Range("N73").Value = -28274.25
$ws.Range("H76").Value = 4405
$ws.Range("J76").Value = 4985
$ws.Range("L76").Value = 4985
$ws.Range("N76").Value = -5615
$ws.Range("H79").Value = 4405
$ws.Range("J79").Value = 4985
$ws.Range("L79").Value = 4985
$ws.Range("N79").Value = -7169
$ws.Range("H86").Value = 2313.6924
$ws.Range("I86").Value = 2175.3333
$ws.Range("K86").Value = 2175.3333
$ws.Range("M86").Value = -1052.3333
$ws.Range("H89").Value = 2313.6924
$ws.Range("I89").Value = 2175.3333
$ws.Range("K89").Value = 10876.6665
$ws.Range("M89").Value = -5260.666499999999
$ws.Range("H100").Value = 2630.75
$ws.Range("I100").Value = 1966
$ws.Range("J100").Value = 4625
$ws.Range("K100").Value = 1966
$ws.Range("L100").Value = 4625
$ws.Range("M100").Value = -1425
$ws.Range("N100").Value = -5707
$ws.Range("H106").Value = 2245
$ws.Range("I106").Value = 1765
$ws.Range("K106").Value = 1765
$ws.Range("M106").Value = -1134
$ws.Range("H132").Value = 3328.879
$ws.Range("I132").Value = 1543.2667
$ws.Range("K132").Value = 4629.800099999999
$ws.Range("M132").Value = -2099.800099999999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2737.2307
$ws.Range("I2").Value = 2235.3635
$ws.Range("K2").Value = 2235.3635
$ws.Range("M2").Value = -2122.3635
$ws.Range("H110").Value = 2215
$ws.Range("I110").Value = 1421
$ws.Range("J110").Value = 4200
$ws.Range("K110").Value = 1421
$ws.Range("L110").Value = 4200
$ws.Range("M110").Value = 624
$ws.Range("N110").Value = -8290
$ws.Range("H116").Value = 2737.2307
$ws.Range("I116").Value = 2235.3635
$ws.Range("K116").Value = 2235.3635
$ws.Range("M116").Value = 58.63650000000007
$ws.Range("H132").Value = 5747.75
$ws.Range("I132").Value = 4330.5
$ws.Range("J132").Value = 9999.5
$ws.Range("K132").Value = 12991.5
$ws.Range("L132").Value = 29998.5
$ws.Range("M132").Value = -10461.5
$ws.Range("N132").Value = -35058.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2737.2307
$ws.Range("I3").Value = 2235.3635
$ws.Range("K3").Value = 2235.3635
$ws.Range("M3").Value = -2121.3635
$ws.Range("H107").Value = 2528.2222
$ws.Range("I107").Value = 2399.3438
$ws.Range("K107").Value = 2399.3438
$ws.Range("M107").Value = -479.3438000000001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2406.4666
$ws.Range("J31").Value = 4305.3125
$ws.Range("L31").Value = 4305.3125
$ws.Range("N31").Value = -4895.3125
$ws.Range("H34").Value = 2406.4666
$ws.Range("J34").Value = 4305.3125
$ws.Range("L34").Value = 4305.3125
$ws.Range("N34").Value = -4709.3125
$ws.Range("H122").Value = 1603.4286
$ws.Range("I122").Value = 1386.2222
$ws.Range("J122").Value = 1994.4
$ws.Range("K122").Value = 4158.6666
$ws.Range("L122").Value = 5983.200000000001
$ws.Range("M122").Value = -1708.6666
$ws.Range("N122").Value = -10883.2
$ws.Range("H132").Value = 2033.76
$ws.Range("I132").Value = 1982.6976
$ws.Range("K132").Value = 5948.0928
$ws.Range("M132").Value = -3418.0928

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 10841
$ws.Range("I113").Value = 14379.625
$ws.Range("J113").Value = 2752.7144
$ws.Range("K113").Value = 14379.625
$ws.Range("L113").Value = 2752.7144
$ws.Range("M113").Value = -12209.625
$ws.Range("N113").Value = -7092.7144
$ws.Range("H123").Value = 56908.5
$ws.Range("J123").Value = 56908.5
$ws.Range("L123").Value = 56908.5
$ws.Range("N123").Value = -61808.5
$ws.Range("H132").Value = 7249.5
$ws.Range("I132").Value = 7570.9287
$ws.Range("K132").Value = 22712.7861
$ws.Range("M132").Value = -20182.7861

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1409.2858
$ws.Range("I16").Value = 1368.5834
$ws.Range("J16").Value = 1653.5
$ws.Range("K16").Value = 1368.5834
$ws.Range("L16").Value = 1653.5
$ws.Range("M16").Value = -1198.5834
$ws.Range("N16").Value = -1993.5
$ws.Range("H22").Value = 1158.591
$ws.Range("I22").Value = 945.75
$ws.Range("J22").Value = 1414
$ws.Range("K22").Value = 945.75
$ws.Range("L22").Value = 1414
$ws.Range("M22").Value = -650.75
$ws.Range("N22").Value = -2004
$ws.Range("H27").Value = 1158.591
$ws.Range("I27").Value = 945.75
$ws.Range("J27").Value = 1414
$ws.Range("K27").Value = 945.75
$ws.Range("L27").Value = 1414
$ws.Range("M27").Value = -838.75
$ws.Range("N27").Value = -1628
$ws.Range("H40").Value = 4055.2083
$ws.Range("I40").Value = 3992.2727
$ws.Range("K40").Value = 3992.2727
$ws.Range("M40").Value = -3856.2727
$ws.Range("H82").Value = 3530.5293
$ws.Range("I82").Value = 3885.7
$ws.Range("J82").Value = 3023.1428
$ws.Range("K82").Value = 3885.7
$ws.Range("L82").Value = 3023.1428
$ws.Range("M82").Value = -3524.7
$ws.Range("N82").Value = -3745.1428
$ws.Range("H85").Value = 3530.5293
$ws.Range("I85").Value = 3885.7
$ws.Range("J85").Value = 3023.1428
$ws.Range("K85").Value = 3885.7
$ws.Range("L85").Value = 3023.1428
$ws.Range("M85").Value = -2637.7
$ws.Range("N85").Value = -5519.1428
$ws.Range("H93").Value = 21828.5
$ws.Range("I93").Value = 952.1667
$ws.Range("K93").Value = 952.1667
$ws.Range("M93").Value = 295.8333
$ws.Range("H122").Value = 4335.769
$ws.Range("I122").Value = 3995.5
$ws.Range("K122").Value = 11986.5
$ws.Range("M122").Value = -9536.5
$ws.Range("H132").Value = 8598.6
$ws.Range("I132").Value = 19999
$ws.Range("J132").Value = 5748.5
$ws.Range("K132").Value = 59997
$ws.Range("L132").Value = 17245.5
$ws.Range("M132").Value = -57467
$ws.Range("N132").Value = -22305.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 101999.5
$ws.Range("I62").Value = 133998.33
$ws.Range("J62").Value = 6003
$ws.Range("K62").Value = 133998.33
$ws.Range("L62").Value = 6003
$ws.Range("M62").Value = -133374.33
$ws.Range("N62").Value = -7251
$ws.Range("H65").Value = 101999.5
$ws.Range("I65").Value = 133998.33
$ws.Range("J65").Value = 6003
$ws.Range("K65").Value = 669991.6499999999
$ws.Range("L65").Value = 30015
$ws.Range("M65").Value = -666871.6499999999
$ws.Range("N65").Value = -36255
$ws.Range("H107").Value = 45507040
$ws.Range("I107").Value = 850.25
$ws.Range("K107").Value = 2550.75
$ws.Range("M107").Value = -630.75
$ws.Range("H132").Value = 4422.2393
$ws.Range("I132").Value = 4508.175
$ws.Range("K132").Value = 13524.525
$ws.Range("M132").Value = -10994.525
